$d = $word.ActiveDocument

# --- helper: wrap a <w:body> fragment in a minimal flat-OPC package so it
#     can be fed to Range.InsertXML (which REPLACES the target range's
#     contents with the parsed fragment). ---
function New-BodyXml {
    param([string]$BodyInner)
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $BodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Work from the bottom of the document upward so earlier paragraph indices
# stay valid while later ones are being rewritten.

# Paragraph 8: "Make new network connections then break old ones"
$p8 = $d.Paragraphs(8)
$xml8 = New-BodyXml(
    '<w:p w14:paraId="1803C610" w14:textId="77777777" w:rsidR="00D27523" w:rsidRDefault="00D27523" w:rsidP="00D27523">' +
      '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
      '<w:r><w:t xml:space="preserve">Make new network connections then break old </w:t></w:r>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:t>ones</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
    '</w:p>'
)
$p8.Range.InsertXML($xml8)

# Paragraph 7: "There are going to be too many graph edges if we run this for any significant number of steps"
$p7 = $d.Paragraphs(7)
$xml7 = New-BodyXml(
    '<w:p w14:paraId="387267C5" w14:textId="77777777" w:rsidR="00D27523" w:rsidRDefault="00D27523" w:rsidP="00D27523">' +
      '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
      '<w:r><w:t xml:space="preserve">There are going to be too many graph edges if we run this for any significant number of </w:t></w:r>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:t>steps</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
    '</w:p>'
)
$p7.Range.InsertXML($xml7)

# Paragraph 6: "Build in a mechanism for pruning graph edges"
$p6 = $d.Paragraphs(6)
$xml6 = New-BodyXml(
    '<w:p w14:paraId="37404A62" w14:textId="77777777" w:rsidR="00D27523" w:rsidRDefault="00D27523" w:rsidP="00D27523">' +
      '<w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' +
      '<w:r><w:t xml:space="preserve">Build in a mechanism for pruning graph </w:t></w:r>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:t>edges</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
    '</w:p>'
)
$p6.Range.InsertXML($xml6)

# Paragraph 3: "Test graph connection dynamic using 25 node toy version of model"
$p3 = $d.Paragraphs(3)
$xml3 = New-BodyXml(
    '<w:p w14:paraId="495D6EEB" w14:textId="6231B5D1" w:rsidR="00D27523" w:rsidRDefault="00D27523" w:rsidP="00D27523">' +
      '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
      '<w:r><w:t xml:space="preserve">Test graph connection dynamic using 25 node toy version of </w:t></w:r>' +
      '<w:proofErr w:type="gramStart"/>' +
      '<w:r><w:t>model</w:t></w:r>' +
      '<w:proofErr w:type="gramEnd"/>' +
    '</w:p>'
)
$p3.Range.InsertXML($xml3)

# Paragraph 1: add the new "Redo everything..." Heading1 + a blank Heading1
# paragraph above the existing "Eight-friends" heading.
$p1 = $d.Paragraphs(1)
$xml1 = New-BodyXml(
    '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Redo everything to be compliant with the new version of agents.jl</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr></w:p>' +
    '<w:p w14:paraId="007D2A2A" w14:textId="2A469762" w:rsidR="00521DD3" w:rsidRDefault="000E4451" w:rsidP="000E4451"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Eight-friends</w:t></w:r></w:p>'
)
$p1.Range.InsertXML($xml1)

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
